# Shorten the five job-title strings in the "Prof" column (D) of NDL6Staff.
# Excel dedups/rebuilds the shared-string table on save and all the sheets
# (GP/GIB/GIM/GOI/GNO) that pull these values via ="NDL6Staff!D#" formulas
# will pick up the new text automatically on recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NDL6Staff")

$ws.Range("D2").Value  = "вед инженер-исследователь"
$ws.Range("D27").Value = "вед инженер-исследователь"

$ws.Range("D7").Value  = "инженер-исследователь I кат"
$ws.Range("D23").Value = "инженер-исследователь I кат"
$ws.Range("D25").Value = "инженер-исследователь I кат"

$ws.Range("D8").Value  = "инженер-исследователь II кат"
$ws.Range("D24").Value = "инженер-исследователь II кат"

$ws.Range("D16").Value = "инженер-электроник I кат"

$ws.Range("D36").Value = "инженер-исследователь без кат"

# Restore the recalculated workbook to match the saved UI state: NDL6Staff
# tab active/selected at D37, GP selection at D5, GIB selection at C11
# (and no longer the active tab).
$wsGP = $wb.Worksheets.Item("GP")
$wsGP.Activate()
$wsGP.Range("D5").Select()

$wsGIB = $wb.Worksheets.Item("GIB")
$wsGIB.Activate()
$wsGIB.Range("C11").Select()

$ws.Activate()
$ws.Range("D37").Select()
